$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dob values for rows 3 and 4 (row 2 keeps its original dob)
$ws.Range("G3").Value = "'2001-04-13"
$ws.Range("G4").Value = "'2001-07-26"

# Update the selection to match the saved view state
$ws.Range("G5").Select()
